{"js": "// Replace the date line and each \"A\u00f7B=\" exercise text with its updated value.\n// Every old value below occurs exactly once in the document body, so a\n// plain case-sensitive search/replace for each pair is unambiguous.\nconst replacements = [\n  [\n    \"2024-04-03 Wednesday\",\n    \"2024-04-04 Thursday\"\n  ],\n  [\n    \"381\u00f78=\",\n    \"500\u00f73=\"\n  ],\n  [\n    \"349\u00f79=\",\n    \"940\u00f73=\"\n  ],\n  [\n    \"431\u00f79=\",\n    \"810\u00f79=\"\n  ],\n  [\n    \"104\u00f77=\",\n    \"656\u00f73=\"\n  ],\n  [\n    \"191\u00f78=\",\n    \"550\u00f72=\"\n  ],\n  [\n    \"379\u00f77=\",\n    \"701\u00f75=\"\n  ],\n  [\n    \"126\u00f75=\",\n    \"203\u00f78=\"\n  ],\n  [\n    \"731\u00f75=\",\n    \"299\u00f76=\"\n  ],\n  [\n    \"586\u00f79=\",\n    \"574\u00f76=\"\n  ],\n  [\n    \"589\u00f75=\",\n    \"861\u00f75=\"\n  ],\n  [\n    \"471\u00f74=\",\n    \"528\u00f76=\"\n  ],\n  [\n    \"408\u00f74=\",\n    \"394\u00f72=\"\n  ],\n  [\n    \"313\u00f73=\",\n    \"950\u00f79=\"\n  ],\n  [\n    \"616\u00f77=\",\n    \"167\u00f79=\"\n  ],\n  [\n    \"165\u00f75=\",\n    \"723\u00f76=\"\n  ],\n  [\n    \"700\u00f74=\",\n    \"268\u00f76=\"\n  ],\n  [\n    \"633\u00f79=\",\n    \"151\u00f74=\"\n  ],\n  [\n    \"592\u00f79=\",\n    \"374\u00f77=\"\n  ],\n  [\n    \"625\u00f79=\",\n    \"375\u00f79=\"\n  ],\n  [\n    \"500\u00f79=\",\n    \"902\u00f74=\"\n  ],\n  [\n    \"581\u00f74=\",\n    \"262\u00f77=\"\n  ],\n  [\n    \"530\u00f72=\",\n    \"281\u00f72=\"\n  ],\n  [\n    \"316\u00f76=\",\n    \"413\u00f75=\"\n  ],\n  [\n    \"144\u00f74=\",\n    \"867\u00f72=\"\n  ],\n  [\n    \"219\u00f74=\",\n    \"396\u00f76=\"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every \"A\u00f7B=\" division exercise to the\n# new values for the day. Each old value occurs exactly once in the\n# document, so a Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-04-03 Wednesday\", \"2024-04-04 Thursday\"),\n    @(\"381\u00f78=\", \"500\u00f73=\"),\n    @(\"349\u00f79=\", \"940\u00f73=\"),\n    @(\"431\u00f79=\", \"810\u00f79=\"),\n    @(\"104\u00f77=\", \"656\u00f73=\"),\n    @(\"191\u00f78=\", \"550\u00f72=\"),\n    @(\"379\u00f77=\", \"701\u00f75=\"),\n    @(\"126\u00f75=\", \"203\u00f78=\"),\n    @(\"731\u00f75=\", \"299\u00f76=\"),\n    @(\"586\u00f79=\", \"574\u00f76=\"),\n    @(\"589\u00f75=\", \"861\u00f75=\"),\n    @(\"471\u00f74=\", \"528\u00f76=\"),\n    @(\"408\u00f74=\", \"394\u00f72=\"),\n    @(\"313\u00f73=\", \"950\u00f79=\"),\n    @(\"616\u00f77=\", \"167\u00f79=\"),\n    @(\"165\u00f75=\", \"723\u00f76=\"),\n    @(\"700\u00f74=\", \"268\u00f76=\"),\n    @(\"633\u00f79=\", \"151\u00f74=\"),\n    @(\"592\u00f79=\", \"374\u00f77=\"),\n    @(\"625\u00f79=\", \"375\u00f79=\"),\n    @(\"500\u00f79=\", \"902\u00f74=\"),\n    @(\"581\u00f74=\", \"262\u00f77=\"),\n    @(\"530\u00f72=\", \"281\u00f72=\"),\n    @(\"316\u00f76=\", \"413\u00f75=\"),\n    @(\"144\u00f74=\", \"867\u00f72=\"),\n    @(\"219\u00f74=\", \"396\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    # wdReplaceAll = 2\n    $find.Execute($null, $true, $null, $null, $null, $null, $true, $null, $null, $null, 2, $true) | Out-Null\n}\n"}
